$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New registered players to append starting at row 14
# Columns: A = Name, C = Phone number (some phone numbers given as text)
$players = @(
    @{ Name = "Fokou Wilfried";   Phone = 691878455;      IsText = $false },
    @{ Name = "Fosso Christian";  Phone = 699054484;      IsText = $false },
    @{ Name = "Happi Steve";      Phone = 694551517;      IsText = $false },
    @{ Name = "Kenfack Dior";     Phone = 656894667;      IsText = $false },
    @{ Name = "Leussi Sahadio";   Phone = 659319642;      IsText = $false },
    @{ Name = "Ngafor Henry";     Phone = 666334062;      IsText = $false },
    @{ Name = "Nzakou Longsen";   Phone = 675723649;      IsText = $false },
    @{ Name = "Pierre";           Phone = 694865719;      IsText = $false },
    @{ Name = "Pola Kouam";       Phone = "6 58 50 39 02"; IsText = $true },
    @{ Name = "Rikam Giovanni";   Phone = "6 55 37 24 22"; IsText = $true },
    @{ Name = "Wouamba Roy";      Phone = 673589923;      IsText = $false }
)

$row = 14
foreach ($p in $players) {
    $ws.Cells.Item($row, 1).Value = $p.Name
    $ws.Cells.Item($row, 3).Value = $p.Phone
    $row++
}

# Scroll the window so row 9 is the top visible row, then select C17
# (mirrors the view state captured in the saved workbook: topLeftCell A9,
# active cell C17)
$ws.Range("C17").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 1
